$wb = $excel.ActiveWorkbook
$day1 = $wb.Worksheets.Item("DAY_01")
$err  = $wb.Worksheets.Item("ERROR_REPORT")

# ---------------------------------------------------------------------------
# DAY_01!row 5 -- new log entry (09.Nov.2016 / CORS task)
# ---------------------------------------------------------------------------

# B5: date-looking text -- must land as literal text (not an Excel date serial)
# and must not disturb the cell's existing style. Stage it on a scratch cell
# formatted as Text, then bring across only the *value* via PasteSpecial so
# the destination keeps its original style index.
$day1.Range("ZZ1").NumberFormat = "@"
$day1.Range("ZZ1").Value = "09.Nov.2016"
$day1.Range("ZZ1").Copy()
$day1.Range("B5").PasteSpecial(-4163)  # xlPasteValues
$day1.Range("ZZ1").Clear()

# C5: plain task description
$day1.Range("C5").Value = "Creating Rest Controller using CORS and access in angular "

# E5: concatenated reference links (as literal text, same trick as B5)
$day1.Range("ZZ1").NumberFormat = "@"
$day1.Range("ZZ1").Value = "https://spring.io/guides/gs/rest-service-cors/                                                                                              https://spring.io/guides/gs/consuming-rest-angularjs/                                                                       "
$day1.Range("ZZ1").Copy()
$day1.Range("E5").PasteSpecial(-4163)
$day1.Range("ZZ1").Clear()

# F5 / G5 / H5: time / errors-encountered flag / numeric column
$day1.Range("F5").Value = "1 Hour"
$day1.Range("G5").Value = "Y"
$day1.Range("H5").Value = 1

# Hyperlink E5 -> the "consuming-rest-angularjs" guide (2nd URL in the cell text)
$day1.Hyperlinks.Add($day1.Range("E5"), "https://spring.io/guides/gs/consuming-rest-angularjs/                                                                       ", "", "", "https://spring.io/guides/gs/consuming-rest-angularjs/                                                                       ")

# Re-apply E5's formatting so it matches the sibling link cells above it
# (E2 carries the same "hyperlink text" style E5 should use).
$day1.Range("E2").Copy()
$day1.Range("E5").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# ERROR_REPORT!row 2 -- CORS error write-up
# ---------------------------------------------------------------------------
$err.Range("A2").Value = 1
$err.Range("B2").WrapText = $true
$err.Range("C2").WrapText = $true
$err.Range("B2").Value = "XMLHttpRequest cannot load http://localhost:8989/arun_online_collaboration/blog/allblogs. No 'Access-Control-Allow-Origin' header is present on the requested resource. Origin 'http://127.0.0.1:8887' is therefore not allowed access."
$err.Range("C2").Value = "Handler method must be annoted with  @CrossOrigin(origins=`"http://127.0.0.1:8887`") and pass the origin (your server where's your second application is running)"
$err.Range("D2").Value = "https://spring.io/guides/gs/rest-service-cors/"

$err.Hyperlinks.Add($err.Range("D2"), "https://spring.io/guides/gs/rest-service-cors/")
# Adding the hyperlink re-styles D2 with a generic "visited link" look;
# restore the sheet's existing link-column style (already used by D3..D8).
$err.Range("D3").Copy()
$err.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# View state: ERROR_REPORT becomes the active/selected tab (D3 selected),
# DAY_01 keeps a remembered selection at C13.
# ---------------------------------------------------------------------------
$day1.Range("C13").Select()
$err.Activate()
$err.Range("D3").Select()

$excel.CutCopyMode = $false
